# Applies the latest price/volume snapshot to the cryptos sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some updated prices are plain decimal numbers (e.g. "142.20"). Excel
# would otherwise silently reinterpret them as floating point numbers and
# drop the significant trailing zero, so force those specific cells to
# Text format first, preserving the exact original string formatting.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "60.981.38"
$ws.Range("E2").Value = "  +0.01%  "
$ws.Range("D3").Value = "3.387.62"
$ws.Range("E3").Value = "  -0.80%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "571.56"
$ws.Range("E5").Value = "  -0.09%  "
$ws.Range("D6").Value = "142.20"
$ws.Range("E6").Value = "  +0.21%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("E8").Value = "  -0.34%  "
$ws.Range("D9").Value = "7.65"
$ws.Range("E9").Value = "  +1.09%  "
$ws.Range("E10").Value = "  -1.38%  "
$ws.Range("E11").Value = "  -0.15%  "
$ws.Range("D12").Value = "3.966.92"
$ws.Range("E12").Value = "  -0.77%  "
$ws.Range("E13").Value = "  +1.92%  "
$ws.Range("D14").Value = "27.70"
$ws.Range("E14").Value = "  -1.83%  "
$ws.Range("D15").Value = "3.415.46"
$ws.Range("E15").Value = "  -0.16%  "
$ws.Range("D16").Value = "0.0000171"
$ws.Range("E16").Value = "  -0.58%  "
$ws.Range("D17").Value = "61.101.29"
$ws.Range("E17").Value = "  -0.01%  "
$ws.Range("D18").Value = "6.10"
$ws.Range("E18").Value = "  -3.44%  "
$ws.Range("D19").Value = "13.65"
$ws.Range("E19").Value = "  -5.32%  "
$ws.Range("D20").Value = "8.95"
$ws.Range("E20").Value = "  -4.44%  "
$ws.Range("D21").Value = "382.16"
$ws.Range("E21").Value = "  -1.97%  "
$ws.Range("D22").Value = "74.82"
$ws.Range("E22").Value = "  +2.79%  "
$ws.Range("E23").Value = "  -2.77%  "
$ws.Range("E24").Value = "  +0.24%  "
$ws.Range("E25").Value = "  -5.48%  "
$ws.Range("D26").Value = "3.520.80"
$ws.Range("E26").Value = "  -0.90%  "
$ws.Range("E27").Value = "  +0.23%  "
$ws.Range("E28").Value = "  +0.01%  "
$ws.Range("D29").Value = "7.32"
$ws.Range("E29").Value = "  -1.50%  "
$ws.Range("E30").Value = "  -0.48%  "
$ws.Range("D31").Value = "7.98"
$ws.Range("E31").Value = "  -2.49%  "
$ws.Range("E32").Value = "  -3.45%  "
$ws.Range("E33").Value = "  -0.01%  "
$ws.Range("D34").Value = "23.33"
$ws.Range("E34").Value = "  -2.35%  "
$ws.Range("D35").Value = "6.97"
$ws.Range("E35").Value = "  -0.54%  "
$ws.Range("D36").Value = "166.15"
$ws.Range("E36").Value = "  -1.05%  "
$ws.Range("D37").Value = "5.03"
$ws.Range("E37").Value = "  -1.94%  "
$ws.Range("D38").Value = "3.417.70"
$ws.Range("E38").Value = "  -0.72%  "
$ws.Range("D39").Value = "1.48"
$ws.Range("E39").Value = "  -4.68%  "
$ws.Range("D40").Value = "0.0770"
$ws.Range("E40").Value = "  -1.71%  "
$ws.Range("D41").Value = "26.92"
$ws.Range("E41").Value = "  -0.58%  "
$ws.Range("E42").Value = "  +0.02%  "
$ws.Range("E43").Value = "  -1.92%  "
$ws.Range("D44").Value = "4.39"
$ws.Range("E44").Value = "  -2.73%  "
$ws.Range("D45").Value = "1.67"
$ws.Range("E45").Value = "  -2.63%  "
$ws.Range("E46").Value = "  -0.65%  "
$ws.Range("D47").Value = "2.452.15"
$ws.Range("E47").Value = "  -6.07%  "
$ws.Range("D48").Value = "23.00"
$ws.Range("E48").Value = "  +0.34%  "
$ws.Range("D49").Value = "6.73"
$ws.Range("E49").Value = "  -3.44%  "
$ws.Range("D50").Value = "0.0266"
$ws.Range("E50").Value = "  +1.73%  "
$ws.Range("D51").Value = "2.14"
$ws.Range("E51").Value = "  +6.55%  "
